$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.130.30'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.138.41'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.36%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '535.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.35'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.24%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.512'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +12.71%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.32'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("E10").Value = '  +7.31%  '

$ws.Range("E11").Value = '  +3.34%  '

$ws.Range("E12").Value = '  +3.42%  '

$ws.Range("E13").Value = '  +1.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000169'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.51%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.239.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.26'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.87%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.136.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.65%  '

$ws.Range("E19").Value = '  +4.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.42%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '377.93'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +9.02%  '

$ws.Range("E22").Value = '  +0.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.88%  '

$ws.Range("E25").Value = '  +4.09%  '

$ws.Range("E26").Value = '  +0.70%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +11.44%  '

$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0885'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.53%  '

$ws.Range("E30").Value = '  +6.00%  '

$ws.Range("E31").Value = '  +1.97%  '

$ws.Range("E32").Value = '  +4.48%  '

$ws.Range("E33").Value = '  +7.56%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.17'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.18%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '161.50'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.18%  '

$ws.Range("E37").Value = '  +10.69%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.62'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.67'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.65%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.638.58'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.46%  '

$ws.Range("E41").Value = '  +6.26%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0677'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.11%  '

$ws.Range("E43").Value = '  +6.54%  '

$ws.Range("E44").Value = '  +1.46%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0273'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.65%  '

$ws.Range("E46").Value = '  -0.05%  '

$ws.Range("E47").Value = '  +4.95%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.101'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +11.75%  '

$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.977'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.75%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.36'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.95%  '

$ws.Range("E51").Value = '  -0.65%  '
